$d = $word.ActiveDocument

$replacements = @(
    @("420÷5=", "370÷3="),
    @("190÷6=", "106÷4="),
    @("320÷9=", "750÷8="),
    @("366÷2=", "513÷5="),
    @("577÷5=", "560÷8="),
    @("450÷4=", "499÷2="),
    @("762÷7=", "994÷2="),
    @("946÷7=", "777÷6="),
    @("387÷4=", "809÷7="),
    @("442÷4=", "962÷6="),
    @("388÷3=", "243÷6="),
    @("723÷5=", "631÷4="),
    @("584÷7=", "969÷8="),
    @("643÷4=", "413÷2="),
    @("183÷9=", "260÷4="),
    @("133÷3=", "773÷8="),
    @("781÷2=", "382÷2="),
    @("981÷6=", "283÷2="),
    @("154÷3=", "949÷2="),
    @("212÷3=", "743÷9="),
    @("525÷9=", "743÷6="),
    @("891÷2=", "550÷4="),
    @("671÷6=", "983÷3="),
    @("118÷2=", "753÷8="),
    @("877÷2=", "440÷8=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
